# Applies the "modified test cases on overdue fix" edit:
#  - Summary sheet: a handful of cells lose the custom "0.00" number-format
#    style (old style index 13) in favour of the plain wrap/valign style
#    already used elsewhere on the sheet (style index 6); two values are
#    bumped by a cent of rounding.
#  - Repayment schedule: the same style cleanup (old indices 9/12/13/14 ->
#    6/7), a handful of rounding tweaks, and a brand new trailing
#    "repayment #7" row (row 9) that zeroes the schedule out.
#  - Transactions: loan/office IDs are renumbered and a couple of Amount /
#    Loan Balance values are corrected to match the schedule; the same
#    style cleanup applies to two cells here too.
#
# Because the unused xf records (numFmtId 15/"no align" and numFmtId 2,
# i.e. the old indices 12/13/14) simply disappear from cellXfs once no
# cell references them any more, we never edit xl/styles.xml directly --
# we just repoint every affected cell at a format that already exists
# (copy/PasteSpecial of formats from a donor cell that already carries the
# desired style) and the unused xf entries are dropped automatically when
# the workbook is saved.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Restyle($ws, $targetRange, $donorCellA1) {
    $ws.Range($donorCellA1).Copy() | Out-Null
    $ws.Range($targetRange).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

# C2 already carries the plain wrap/valign style (old index 6) -- use it
# as the donor for every cell that needs to drop its "0.00"/"#,##0" look.
Restyle $summary "B2" "C2"
Restyle $summary "F2" "C2"
Restyle $summary "A3" "C2"
Restyle $summary "B3" "C2"
Restyle $summary "E3" "C2"
Restyle $summary "F3" "C2"

$summary.Range("A3").Value = 106.04
$summary.Range("E3").Value = 89.05

# Final on-sheet selection left by the author (a whole-row block below
# the data, top-left cell A7 is naturally the active cell here).
$summary.Range("A7:XFD15").Select() | Out-Null

# ---------------------------------------------------------------------
# Repayment schedule
# ---------------------------------------------------------------------
$sched = $wb.Worksheets.Item("Repayment schedule")

# Bulk style cleanup: everything below goes to the plain wrap/valign
# style already used by A2 (old index 6).
Restyle $sched "F3:F8" "A2"
Restyle $sched "H3:H8" "A2"
Restyle $sched "K3:K8" "A2"
Restyle $sched "M3:M8" "A2"
Restyle $sched "P4:P8" "A2"
Restyle $sched "G7:G8" "A2"
Restyle $sched "L3"    "A2"

# D3 moves from the old "date, no align" style to the date+wrap style
# already used by C2:C8 (old index 7).
Restyle $sched "D3" "C2"

# E3 moves from the plain wrap/valign style to the italic-font wrap/valign
# style already used on the Transactions sheet (old index 10).
$trans = $wb.Worksheets.Item("Transactions")
$trans.Range("K2").Copy() | Out-Null
$sched.Range("E3").PasteSpecial($xlPasteFormats) | Out-Null

# Value corrections (rounding fixes).
$sched.Range("L3").Value = 850.32
$sched.Range("G4").Value = 3333.34
$sched.Range("G5").Value = 2500.0100000000002
$sched.Range("G6").Value = 1666.68
$sched.Range("G7").Value = 833.35
$sched.Range("G8").Value = 0.02

# New trailing row describing a 7th, fully-settled repayment.
$sched.Range("A9").Value = 7
$sched.Range("B9").Value = 31
$sched.Range("C9").Value = 42217
$sched.Range("D9").Value = ""
$sched.Range("E9").Value = ""
$sched.Range("F9").Value = 0.02
$sched.Range("G9").Value = 0
$sched.Range("H9").Value = 0
$sched.Range("I9").Value = 0
$sched.Range("J9").Value = 0
$sched.Range("K9").Value = 0.02
$sched.Range("L9").Value = 0
$sched.Range("M9").Value = 0
$sched.Range("N9").Value = 0
$sched.Range("O9").Value = 0
$sched.Range("P9").Value = 0.02

# Row 9 uses the plain wrap/valign style throughout, except C9 which is a
# date and gets the date+wrap style (same donors as above).
Restyle $sched "A9:B9" "A2"
Restyle $sched "D9:P9" "A2"
Restyle $sched "C9" "C2"

# Final on-sheet selection left by the author (the blank row right below
# the newly appended row 9).
$sched.Range("A10:XFD10").Select() | Out-Null

# ---------------------------------------------------------------------
# Transactions
# ---------------------------------------------------------------------
Restyle $trans "F2" "H2"
Restyle $trans "G2" "H2"

$trans.Range("A2").Value = 96
$trans.Range("E2").Value = 850.32
$trans.Range("J2").Value = 4166.67
$trans.Range("A3").Value = 94

# Final on-sheet selection left by the author: the whole A2:XFD4 block
# (active cell anchored at A4 in the source file -- a product of a
# shift-click the user made while reviewing; Select() lands the active
# cell on the range's first cell, which is the closest reproducible
# approximation here).
$trans.Range("A2:XFD4").Select() | Out-Null
